# Fruta / hortaliza, semanal
# Insert three new weekly price rows for "Femacal de La Calera" / Ciruela / Black Amber
# (Provincia de San Felipe de Aconcagua) above the existing row 122, pushing the
# existing rows 122-132 down to 125-135.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 122:124 - shifts old rows 122-132 down to 125-135
$ws.Rows("122:124").Insert()

# Common columns for the block (constant across this product's rows)
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100103
$producto  = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria = "Ciruela"
$variedad  = "Black Amber"
$unidad    = "$/caja 15 kilos empedrada"
$origen    = "Provincia de San Felipe de Aconcagua"
$fecha     = 44578

# Row 122: Especial
$ws.Cells.Item(122, 1).Value  = $mercadoId
$ws.Cells.Item(122, 2).Value  = $mercado
$ws.Cells.Item(122, 3).Value  = $region
$ws.Cells.Item(122, 4).Value  = $fecha
$ws.Cells.Item(122, 5).Value  = $codreg
$ws.Cells.Item(122, 6).Value  = $tipo
$ws.Cells.Item(122, 7).Value  = $productoId
$ws.Cells.Item(122, 8).Value  = $producto
$ws.Cells.Item(122, 9).Value  = $categoriaId
$ws.Cells.Item(122, 10).Value = $categoria
$ws.Cells.Item(122, 11).Value = $variedad
$ws.Cells.Item(122, 12).Value = "Especial"
$ws.Cells.Item(122, 13).Value = 85
$ws.Cells.Item(122, 14).Value = 13000
$ws.Cells.Item(122, 15).Value = 13000
$ws.Cells.Item(122, 16).Value = 13000
$ws.Cells.Item(122, 17).Value = $unidad
$ws.Cells.Item(122, 18).Value = $origen
$ws.Cells.Item(122, 19).Value = 867
$ws.Cells.Item(122, 20).Value = 15

# Row 123: Primera
$ws.Cells.Item(123, 1).Value  = $mercadoId
$ws.Cells.Item(123, 2).Value  = $mercado
$ws.Cells.Item(123, 3).Value  = $region
$ws.Cells.Item(123, 4).Value  = $fecha
$ws.Cells.Item(123, 5).Value  = $codreg
$ws.Cells.Item(123, 6).Value  = $tipo
$ws.Cells.Item(123, 7).Value  = $productoId
$ws.Cells.Item(123, 8).Value  = $producto
$ws.Cells.Item(123, 9).Value  = $categoriaId
$ws.Cells.Item(123, 10).Value = $categoria
$ws.Cells.Item(123, 11).Value = $variedad
$ws.Cells.Item(123, 12).Value = "Primera"
$ws.Cells.Item(123, 13).Value = 80
$ws.Cells.Item(123, 14).Value = 12000
$ws.Cells.Item(123, 15).Value = 12000
$ws.Cells.Item(123, 16).Value = 12000
$ws.Cells.Item(123, 17).Value = $unidad
$ws.Cells.Item(123, 18).Value = $origen
$ws.Cells.Item(123, 19).Value = 800
$ws.Cells.Item(123, 20).Value = 15

# Row 124: Segunda
$ws.Cells.Item(124, 1).Value  = $mercadoId
$ws.Cells.Item(124, 2).Value  = $mercado
$ws.Cells.Item(124, 3).Value  = $region
$ws.Cells.Item(124, 4).Value  = $fecha
$ws.Cells.Item(124, 5).Value  = $codreg
$ws.Cells.Item(124, 6).Value  = $tipo
$ws.Cells.Item(124, 7).Value  = $productoId
$ws.Cells.Item(124, 8).Value  = $producto
$ws.Cells.Item(124, 9).Value  = $categoriaId
$ws.Cells.Item(124, 10).Value = $categoria
$ws.Cells.Item(124, 11).Value = $variedad
$ws.Cells.Item(124, 12).Value = "Segunda"
$ws.Cells.Item(124, 13).Value = 75
$ws.Cells.Item(124, 14).Value = 10000
$ws.Cells.Item(124, 15).Value = 10000
$ws.Cells.Item(124, 16).Value = 10000
$ws.Cells.Item(124, 17).Value = $unidad
$ws.Cells.Item(124, 18).Value = $origen
$ws.Cells.Item(124, 19).Value = 667
$ws.Cells.Item(124, 20).Value = 15
